$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf1"
$ws.Cells.Item(2,3).Value = "Fgfr3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.8775636666666666
$ws.Cells.Item(2,8).Value = 2.632691
$ws.Cells.Item(2,9).Value = 0.1887436506618166
$ws.Cells.Item(2,10).Value = 0.2083714858314108
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.330840333333333
$ws.Cells.Item(2,14).Value = 6.992521
$ws.Cells.Item(2,15).Value = 0.6715345129768794
$ws.Cells.Item(2,16).Value = 0.7003397275969581
$ws.Cells.Item(2,17).Value = 2.045460789334555
$ws.Cells.Item(2,18).Value = 18.409147104011
$ws.Cells.Item(2,19).Value = 0.1267478755246613
$ws.Cells.Item(2,20).Value = 0.1459308296261436

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf1"
$ws.Cells.Item(3,3).Value = "Fgfr3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.8775636666666666
$ws.Cells.Item(3,8).Value = 2.632691
$ws.Cells.Item(3,9).Value = 0.1887436506618166
$ws.Cells.Item(3,10).Value = 0.2083714858314108
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.657666
$ws.Cells.Item(3,14).Value = 1.972998
$ws.Cells.Item(3,15).Value = 0.1894790521235985
$ws.Cells.Item(3,16).Value = 0.1976066831789769
$ws.Cells.Item(3,17).Value = 0.5771437864019999
$ws.Cells.Item(3,18).Value = 5.194294077618
$ws.Cells.Item(3,19).Value = 0.03576296802174861
$ws.Cells.Item(3,20).Value = 0.04117559818422027

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf1"
$ws.Cells.Item(4,3).Value = "Fgfr3"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.8775636666666666
$ws.Cells.Item(4,8).Value = 2.632691
$ws.Cells.Item(4,9).Value = 0.1887436506618166
$ws.Cells.Item(4,10).Value = 0.2083714858314108
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.42828
$ws.Cells.Item(4,14).Value = 0.85656
$ws.Cells.Item(4,15).Value = 0.123391035029171
$ws.Cells.Item(4,16).Value = 0.08578923067523865
$ws.Cells.Item(4,17).Value = 0.37584296716
$ws.Cells.Item(4,18).Value = 2.25505780296
$ws.Cells.Item(4,19).Value = 0.02328927441034582
$ws.Cells.Item(4,20).Value = 0.01787602946413312

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Fgf1"
$ws.Cells.Item(5,3).Value = "Fgfr3"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.8775636666666666
$ws.Cells.Item(5,8).Value = 2.632691
$ws.Cells.Item(5,9).Value = 0.1887436506618166
$ws.Cells.Item(5,10).Value = 0.2083714858314108
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.05413033333333334
$ws.Cells.Item(5,14).Value = 0.162391
$ws.Cells.Item(5,15).Value = 0.01559539987035126
$ws.Cells.Item(5,16).Value = 0.01626435854882633
$ws.Cells.Item(5,17).Value = 0.04750281379788889
$ws.Cells.Item(5,18).Value = 0.427525324181
$ws.Cells.Item(5,19).Value = 0.002943532705060917
$ws.Cells.Item(5,20).Value = 0.00338902855691375

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf1"
$ws.Cells.Item(6,3).Value = "Fgfr3"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.458038666666667
$ws.Cells.Item(6,8).Value = 7.374116000000001
$ws.Cells.Item(6,9).Value = 0.5286672739959656
$ws.Cells.Item(6,10).Value = 0.5836444564186148
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.330840333333333
$ws.Cells.Item(6,14).Value = 6.992521
$ws.Cells.Item(6,15).Value = 0.6715345129768794
$ws.Cells.Item(6,16).Value = 0.7003397275969581
$ws.Cells.Item(6,17).Value = 5.729295665159555
$ws.Cells.Item(6,18).Value = 51.56366098643601
$ws.Cells.Item(6,19).Value = 0.3550183203696952
$ws.Cells.Item(6,20).Value = 0.4087493996216874

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf1"
$ws.Cells.Item(7,3).Value = "Fgfr3"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.458038666666667
$ws.Cells.Item(7,8).Value = 7.374116000000001
$ws.Cells.Item(7,9).Value = 0.5286672739959656
$ws.Cells.Item(7,10).Value = 0.5836444564186148
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.657666
$ws.Cells.Item(7,14).Value = 1.972998
$ws.Cells.Item(7,15).Value = 0.1894790521235985
$ws.Cells.Item(7,16).Value = 0.1976066831789769
$ws.Cells.Item(7,17).Value = 1.616568457752
$ws.Cells.Item(7,18).Value = 14.549116119768
$ws.Cells.Item(7,19).Value = 0.1001713739655223
$ws.Cells.Item(7,20).Value = 0.1153320451886794

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Fgf1"
$ws.Cells.Item(8,3).Value = "Fgfr3"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.458038666666667
$ws.Cells.Item(8,8).Value = 7.374116000000001
$ws.Cells.Item(8,9).Value = 0.5286672739959656
$ws.Cells.Item(8,10).Value = 0.5836444564186148
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.42828
$ws.Cells.Item(8,14).Value = 0.85656
$ws.Cells.Item(8,15).Value = 0.123391035029171
$ws.Cells.Item(8,16).Value = 0.08578923067523865
$ws.Cells.Item(8,17).Value = 1.05272880016
$ws.Cells.Item(8,18).Value = 6.316372800960001
$ws.Cells.Item(8,19).Value = 0.06523280212441253
$ws.Cells.Item(8,20).Value = 0.05007040890402082

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Fgf1"
$ws.Cells.Item(9,3).Value = "Fgfr3"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.458038666666667
$ws.Cells.Item(9,8).Value = 7.374116000000001
$ws.Cells.Item(9,9).Value = 0.5286672739959656
$ws.Cells.Item(9,10).Value = 0.5836444564186148
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.05413033333333334
$ws.Cells.Item(9,14).Value = 0.162391
$ws.Cells.Item(9,15).Value = 0.01559539987035126
$ws.Cells.Item(9,16).Value = 0.01626435854882633
$ws.Cells.Item(9,17).Value = 0.1330544523728889
$ws.Cells.Item(9,18).Value = 1.197490071356
$ws.Cells.Item(9,19).Value = 0.008244777536335634
$ws.Cells.Item(9,20).Value = 0.009492602704227193

$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Fgf1"
$ws.Cells.Item(10,3).Value = "Fgfr3"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.313898
$ws.Cells.Item(10,8).Value = 2.627796
$ws.Cells.Item(10,9).Value = 0.2825890753422177
$ws.Cells.Item(10,10).Value = 0.2079840577499744
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.330840333333333
$ws.Cells.Item(10,14).Value = 6.992521
$ws.Cells.Item(10,15).Value = 0.6715345129768794
$ws.Cells.Item(10,16).Value = 0.7003397275969581
$ws.Cells.Item(10,17).Value = 3.062486452286
$ws.Cells.Item(10,18).Value = 18.374918713716
$ws.Cells.Item(10,19).Value = 0.1897683170825228
$ws.Cells.Item(10,20).Value = 0.1456594983491271

$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Fgf1"
$ws.Cells.Item(11,3).Value = "Fgfr3"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.313898
$ws.Cells.Item(11,8).Value = 2.627796
$ws.Cells.Item(11,9).Value = 0.2825890753422177
$ws.Cells.Item(11,10).Value = 0.2079840577499744
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.657666
$ws.Cells.Item(11,14).Value = 1.972998
$ws.Cells.Item(11,15).Value = 0.1894790521235985
$ws.Cells.Item(11,16).Value = 0.1976066831789769
$ws.Cells.Item(11,17).Value = 0.864106042068
$ws.Cells.Item(11,18).Value = 5.184636252408001
$ws.Cells.Item(11,19).Value = 0.05354471013632758
$ws.Cells.Item(11,20).Value = 0.04109903980607724

$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Fgf1"
$ws.Cells.Item(12,3).Value = "Fgfr3"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.313898
$ws.Cells.Item(12,8).Value = 2.627796
$ws.Cells.Item(12,9).Value = 0.2825890753422177
$ws.Cells.Item(12,10).Value = 0.2079840577499744
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.42828
$ws.Cells.Item(12,14).Value = 0.85656
$ws.Cells.Item(12,15).Value = 0.123391035029171
$ws.Cells.Item(12,16).Value = 0.08578923067523865
$ws.Cells.Item(12,17).Value = 0.56271623544
$ws.Cells.Item(12,18).Value = 2.25086494176
$ws.Cells.Item(12,19).Value = 0.03486895849441263
$ws.Cells.Item(12,20).Value = 0.01784279230708471

$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Fgf1"
$ws.Cells.Item(13,3).Value = "Fgfr3"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.313898
$ws.Cells.Item(13,8).Value = 2.627796
$ws.Cells.Item(13,9).Value = 0.2825890753422177
$ws.Cells.Item(13,10).Value = 0.2079840577499744
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.05413033333333334
$ws.Cells.Item(13,14).Value = 0.162391
$ws.Cells.Item(13,15).Value = 0.01559539987035126
$ws.Cells.Item(13,16).Value = 0.01626435854882633
$ws.Cells.Item(13,17).Value = 0.071121736706
$ws.Cells.Item(13,18).Value = 0.426730420236
$ws.Cells.Item(13,19).Value = 0.004407089628954704
$ws.Cells.Item(13,20).Value = 0.003382727287685385
